# Dagbaekur Ivan Gusti Saevar uddfaerdar
# Fill in missing time-tracking entries (in minutes) for "Vika 6", "Vika 7"
# and "Vika 8" on Sheet1, then update the yearly summary (row 62) to sum
# the weekly totals directly and add a "hours" helper cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Vika 6 (rows 26-32): Hönnun (row 28) ---
$ws.Range("C28").Value = 60

# --- Vika 7 (rows 36-42): Hönnun (38), Rannsóknir (39), Forritun (40) ---
$ws.Range("C38").Value = 120
$ws.Range("D38").Value = 60

$ws.Range("C39").Value = 180
$ws.Range("D39").Value = 180
$ws.Range("E39").Value = 240

$ws.Range("D40").Value = 120
$ws.Range("E40").Value = 120

# --- Vika 8 (rows 46-52): Kröfulýsing (46), Hönnun (48), Rannsóknir (49),
#     Forritun (50), Prófanir (51) ---
$ws.Range("D46").Value = 60

$ws.Range("C48").Value = 60
$ws.Range("D48").Value = 60
$ws.Range("E48").Value = 60

$ws.Range("C49").Value = 120
$ws.Range("D49").Value = 180
$ws.Range("E49").Value = 120

$ws.Range("D50").Value = 60

$ws.Range("D51").Value = 60
$ws.Range("E51").Value = 120

# --- Samantekt (summary): total time now sums the weekly totals directly,
#     plus a new "hours" column next to it ---
$ws.Range("D62").Formula = "=SUM(J12,J22,J32,J42,J52)"
$ws.Range("F62").Formula = "=D62/60"

# Leave the view pointed near the bottom of the sheet, where the new data was entered.
$ws.Range("F59").Select()
